# Rename the sheet from "Property1" to "DataNode" — unifying the
# DataNode / DataTable / Entity naming convention referenced in the
# commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Match the author's on-disk cursor/selection position recorded in the
# saved view state (pane stays frozen at row 9; active cell moves to D36).
$ws.Range("D36").Select() | Out-Null
